# Apply the edit described in the diff:
# Two new rows of data were inserted into the "Pera" price sheet between
# the rows that used to be 419 and 420 (1-indexed), pushing the former
# rows 420-493 down to 422-495. The sheet dimension grows from
# A1:T493 to A1:T495.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 420; this shifts the existing
# rows 420:493 down to 422:495 and copies formatting (incl. the date
# number format in column D) from the row above, just like Excel's
# native "Insert Sheet Rows" command.
$ws.Rows("420:421").Insert()

# --- Populate the first new row (420) ---
$ws.Range("A420").Value = 10
$ws.Range("B420").Value = "Vega Modelo de Temuco"
$ws.Range("C420").Value = "La Araucanía"
$ws.Range("D420").Value = 44505
$ws.Range("E420").Value = 9
$ws.Range("F420").Value = "Fruta"
$ws.Range("G420").Value = 100104
$ws.Range("H420").Value = "Frutos de pepita"
$ws.Range("I420").Value = 100104005
$ws.Range("J420").Value = "Pera"
$ws.Range("K420").Value = "Packham's Triumph"
$ws.Range("L420").Value = "Calibre 80"
$ws.Range("M420").Value = 55
$ws.Range("N420").Value = 24000
$ws.Range("O420").Value = 24000
$ws.Range("P420").Value = 24000
$ws.Range("Q420").Value = "`$/caja 18 kilos embalada"
$ws.Range("R420").Value = "Región de O'Higgins"
$ws.Range("S420").Value = 1333
$ws.Range("T420").Value = 18

# --- Populate the second new row (421) ---
$ws.Range("A421").Value = 10
$ws.Range("B421").Value = "Vega Modelo de Temuco"
$ws.Range("C421").Value = "La Araucanía"
$ws.Range("D421").Value = 44505
$ws.Range("E421").Value = 9
$ws.Range("F421").Value = "Fruta"
$ws.Range("G421").Value = 100104
$ws.Range("H421").Value = "Frutos de pepita"
$ws.Range("I421").Value = 100104005
$ws.Range("J421").Value = "Pera"
$ws.Range("K421").Value = "Packham's Triumph"
$ws.Range("L421").Value = "Primera"
$ws.Range("M421").Value = 110
$ws.Range("N421").Value = 14000
$ws.Range("O421").Value = 14000
$ws.Range("P421").Value = 14000
$ws.Range("Q421").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R421").Value = "Región de O'Higgins"
$ws.Range("S421").Value = 778
$ws.Range("T421").Value = 18
